# Generate Report for Handback
#
# This localization-status report is regenerated whenever a handback is
# produced.  The handback for the 1e5cb89d...md document completed, so:
#   - the Overview sheet's per-locale status cells move from
#     "Ready for handoff" to "Handed back: in sync with en-US"
#   - each locale sheet (zh-cn / de-de) gets its "Latest Target File" (I)
#     and "Latest Handback File" (J) columns populated, with I becoming a
#     hyperlink to the source markdown doc (matching column A's link)
#   - the "Latest Handback DateTime" (K) timestamp is refreshed
#   - a few report columns are widened to fit the newly-populated values

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

$docDisplay = "1e5cb89d-5659-48d6-9c61-93610b0ac39f.md"
$docUrl     = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/8625efa7ab2e2eeb4bd64c1f5cdb3813ad9c8e52/e2e/1e5cb89d-5659-48d6-9c61-93610b0ac39f.md"

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: status columns for zh-cn (E2) and de-de (F2) ---
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus

# Widen the Overview status columns to fit the longer text.
$overview.Range("E1").ColumnWidth = 29.9777047293527
$overview.Range("F1").ColumnWidth = 29.9777047293527

# --- zh-cn sheet ---
$zhcn.Range("C2").Value = $newStatus

# Latest Target File (I2) -> hyperlink to the source document.
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $docUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $docDisplay) | Out-Null

# Latest Handback File (J2) -> generated xliff handback file name.
$zhcn.Range("J2").Value = "1e5cb89d-5659-48d6-9c61-93610b0ac39f.09aa81012602d30b190cf72cc7991add4a918cf5.zh-cn.xlf"

# Latest Handback DateTime (K2) refreshed.
$zhcn.Range("K2").Value = "2016-09-07 01:18:31"

# Widen columns C, I, J to fit new content.
$zhcn.Range("C1").ColumnWidth = 29.9777047293527
$zhcn.Range("I1").ColumnWidth = 40
$zhcn.Range("J1").ColumnWidth = 40

# --- de-de sheet ---
$dede.Range("C2").Value = $newStatus

# Latest Target File (I2) -> hyperlink to the source document.
$dede.Hyperlinks.Add($dede.Range("I2"), $docUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $docDisplay) | Out-Null

# Latest Handback File (J2) -> generated xliff handback file name.
$dede.Range("J2").Value = "1e5cb89d-5659-48d6-9c61-93610b0ac39f.09aa81012602d30b190cf72cc7991add4a918cf5.de-de.xlf"

# Latest Handback DateTime (K2) refreshed (later than zh-cn's).
$dede.Range("K2").Value = "2016-09-07 01:18:40"

# Widen columns C, I, J to fit new content.
$dede.Range("C1").ColumnWidth = 29.9777047293527
$dede.Range("I1").ColumnWidth = 40
$dede.Range("J1").ColumnWidth = 40
